$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellB64 {
    param($Sheet, $Addr, $B64)
    $bytes = [Convert]::FromBase64String($B64)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $Sheet.Range($Addr).Value = $text
}

# Row 2
Set-CellB64 $ws "C2" "VW5rbm93biBUaXRsZQ=="
Set-CellB64 $ws "D2" "VW5rbm93biBBYnN0cmFjdA=="
Set-CellB64 $ws "E2" "W10="
Set-CellB64 $ws "F2" "bm90IGZvdW5k"
Set-CellB64 $ws "G2" "Ti9B"
Set-CellB64 $ws "H2" "MTk3MC0wMS0wMQ=="
Set-CellB64 $ws "I2" ""

# Row 3
Set-CellB64 $ws "D3" "QmFja2dyb3VuZDogU2VsZi1zdGlnbWEgb2NjdXJzIHdoZW4gcGVvcGxlIHdpdGggbWVudGFsIGlsbG5lc3NlcyBpbnRlcm5hbGl6ZSBuZWdhdGl2ZSBzdGVyZW90eXBlcyBhbmQgcHJlanVkaWNlcyBhYm91dCB0aGVpciBjb25kaXRpb24uCiBJdCBjYW4gcmVkdWNlIGhlbHAtc2Vla2luZyBiZWhhdmlvdXIgYW5kIHRyZWF0bWVudCBhZGhlcmVuY2UuCiBUaGUgZWZmZWN0aXZlbmVzcyBvZiBpbnRlcnZlbnRpb25zIGFpbWVkIGF0IHJlZHVjaW5nIHNlbGYtc3RpZ21hIGluIHBlb3BsZSB3aXRoIG1lbnRhbCBpbGxuZXNzIGlzIHN5c3RlbWF0aWNhbGx5IHJldmlld2VkLgogUmVzdWx0cyBhcmUgZGlzY3Vzc2VkIGluIHRoZSBjb250ZXh0IG9mIGEgbG9naWMgbW9kZWwgb2YgdGhlIGJyb2FkZXIgc29jaWFsIGNvbnRleHQgb2YgbWVudGFsIGlsbG5lc3Mgc3RpZ21hLgog"
Set-CellB64 $ws "E3" "W1JvbGFuZCBCcmlhbiVCw7xjaHRlciVOVUxMJTAsIE1lbGFuaWUlTWVzc2VyJU5VTEwlMV0="
Set-CellB64 $ws "I3" ""
Set-CellB64 $ws "J3" "R2VybWFuIE1lZGljYWwgU2NpZW5jZSBHTVMgUHVibGlzaGluZyBIb3VzZQ=="

# Row 4
Set-CellB64 $ws "C4" "VW5rbm93biBUaXRsZQ=="
Set-CellB64 $ws "E4" "W10="
Set-CellB64 $ws "F4" "bm90IGZvdW5k"
Set-CellB64 $ws "G4" "Ti9B"
Set-CellB64 $ws "H4" "MTk3MC0wMS0wMQ=="
Set-CellB64 $ws "J4" ""

# Row 5
Set-CellB64 $ws "E5" "W0dhcnVtbWEgVG9sdSVGZXlpc3NhJU5VTEwlMCwgQ3JhaWclTG9ja3dvb2QlTlVMTCUyLCBDcmFpZyVMb2Nrd29vZCVOVUxMJTAsIE1pcmt1emllJVdvbGRpZSVOVUxMJTEsIFphY2hhcnklTXVubiVOVUxMJTEsIEplcm9tZSBBLiVTaW5naCVOVUxMJTIsIEplcm9tZSBBLiVTaW5naCVOVUxMJTAsIE5VTEwlTlVMTCVOVUxMJTAsIE5VTEwlTlVMTCVOVUxMJTBd"
Set-CellB64 $ws "I5" ""
Set-CellB64 $ws "J5" "UHVibGljIExpYnJhcnkgb2YgU2NpZW5jZQ=="

# Row 6
Set-CellB64 $ws "E6" "W1NhYmluZSBFLiVIYW5pc2NoJWhhbmlzY2guc2FiaW5lQHNpZW1lbnMuY29tJTAsIENvbmFsIEQuJVR3b21leSVjLnR3b21leUBzb3Rvbi5hYy51ayUxLCBBbmRyZXcgQy4gSC4lU3pldG8lYXN6ZXRvQHVjYWxnYXJ5LmNhJTEsIFVscmljaCBXLiVCaXJuZXIldWxyaWNoLmJpcm5lckBzaWVtZW5zLmNvbSUxLCBEZW5uaXMlTm93YWslZGVubmlzLm5vd2FrQG1lZC51bmktbXVlbmNoZW4uZGUlMSwgQ2FybGElU2FiYXJpZWdvJUNhcmxhLnNhYmFyaWVnb0BtZWQudW5pLW11ZW5jaGVuLmRlJTFd"
Set-CellB64 $ws "I6" ""
Set-CellB64 $ws "J6" "QmlvTWVkIENlbnRyYWw="

# Row 7
Set-CellB64 $ws "E7" "W0UuJUhlaW0lTlVMTCUwLCBCLiBBLiVLb2hydCVOVUxMJTIsIEIuIEEuJUtvaHJ0JU5VTEwlMCwgTS4lS29zY2hvcmtlJU5VTEwlMSwgTS4lTWlsZW5vdmElTlVMTCUxLCBHLiVUaG9ybmljcm9mdCVOVUxMJTFd"
Set-CellB64 $ws "I7" ""
Set-CellB64 $ws "J7" "Q2FtYnJpZGdlIFVuaXZlcnNpdHkgUHJlc3M="

# Row 8
Set-CellB64 $ws "C8" "VW5rbm93biBUaXRsZQ=="
Set-CellB64 $ws "D8" "VW5rbm93biBBYnN0cmFjdA=="
Set-CellB64 $ws "E8" "W10="
Set-CellB64 $ws "F8" "bm90IGZvdW5k"
Set-CellB64 $ws "G8" "Ti9B"
Set-CellB64 $ws "H8" "MTk3MC0wMS0wMQ=="
Set-CellB64 $ws "I8" ""

# Row 9
Set-CellB64 $ws "C9" "VW5rbm93biBUaXRsZQ=="
Set-CellB64 $ws "D9" "VW5rbm93biBBYnN0cmFjdA=="
Set-CellB64 $ws "E9" "W10="
Set-CellB64 $ws "F9" "bm90IGZvdW5k"
Set-CellB64 $ws "G9" "Ti9B"
Set-CellB64 $ws "H9" "MTk3MC0wMS0wMQ=="
Set-CellB64 $ws "I9" ""

# Row 10
Set-CellB64 $ws "E10" "W01hcmlqYSVQYW50ZWxpYyVOVUxMJTAsIEphbmluYSBJJVN0ZWluZXJ0JU5VTEwlMiwgSmFuaW5hIEklU3RlaW5lcnQlTlVMTCUwLCBKYXklUGFyayVOVUxMJTEsIFNoYXVuJU1lbGxvcnMlTlVMTCUxLCBGdW5nYWklTXVyYXUlTlVMTCUxXQ=="
Set-CellB64 $ws "I10" ""
Set-CellB64 $ws "J10" "Qk1KIFB1Ymxpc2hpbmcgR3JvdXA="

# Row 11
Set-CellB64 $ws "C11" "VW5rbm93biBUaXRsZQ=="
Set-CellB64 $ws "E11" "W10="
Set-CellB64 $ws "F11" "bm90IGZvdW5k"
Set-CellB64 $ws "G11" "Ti9B"
Set-CellB64 $ws "H11" "MTk3MC0wMS0wMQ=="
Set-CellB64 $ws "J11" ""

# Row 12
Set-CellB64 $ws "C12" "VW5rbm93biBUaXRsZQ=="
Set-CellB64 $ws "F12" "bm90IGZvdW5k"
Set-CellB64 $ws "G12" "Ti9B"
Set-CellB64 $ws "I12" ""
